$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Add extra Larman citation detail to the two "reading direction
#    (Larman, 2004)" occurrences (review + suggestions sections).
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "reading direction (Larman, 2004).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "reading direction (Larman, chapter 9, 9.12 figure, 2004).",
    2) | Out-Null

# ------------------------------------------------------------------
# 2) Req. 1 paragraph: extend the sentence with the extra reasoning
#    and citation.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "could therefore be put in a separate class.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "could therefore be put in a separate class for the two to make use of (Larman, chapter 9, figure 9.9, 2004).",
    2) | Out-Null

# ------------------------------------------------------------------
# 3) Suggestions list: "... for the two to make use of (Larman, 2004)."
#    gets the extra chapter/figure detail too.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "for the two to make use of (Larman, 2004).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "for the two to make use of (Larman, chapter 9, figure 9.9, 2004).",
    2) | Out-Null

# ------------------------------------------------------------------
# 4) Both "[0..1]." occurrences (Picture attribute) gain a citation
#    right before the final period.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    " [0..1].",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " [0..1] (Larman, chapter 9, figure 9.20, 2004).",
    2) | Out-Null

# ------------------------------------------------------------------
# 5) Remove two superfluous blank paragraphs from the run of blank
#    paragraphs right before the "Review Content" heading.
# ------------------------------------------------------------------
$idx = 1
while ($idx -le $d.Paragraphs.Count) {
    if ($d.Paragraphs($idx).Range.Text.Trim() -eq "Review Content") {
        break
    }
    $idx = $idx + 1
}

$d.Paragraphs($idx - 1).Range.Delete()
$d.Paragraphs($idx - 2).Range.Delete()
